# Student-answer workbook: convert the "student id" column (A1:A80) from the
# padded 7001-7080 series down to a plain 1-80 integer id, and tidy up the
# worksheet view so it opens scrolled to the top with A1:A80 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data fix: A1:A80 held 7001..7080, now should just be the row number ---
for ($row = 1; $row -le 80; $row++) {
    $ws.Cells.Item($row, 1).Value = $row
}

# --- View fix: scroll back to the top-left and select A1:A80 -------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1:A80").Select()
